$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (21 and 22) for the new NHR.Consensus keys, and fix the
# existing F8 translation, in the exact order needed to reproduce the
# resulting shared-string table ordering.
$ws.Range("E22").Value = "No more mutual reject."
$ws.Range("C21").Value = "NHR.Consensus"
$ws.Range("C22").Value = "NHR.Consensus_Tip"
$ws.Range("F8").Value = "궁합 임계값:"
$ws.Range("F22").Value = "더 이상 상호 거부가 없습니다."
$ws.Range("E21").Value = "Consensus romance always succeed"
$ws.Range("F21").Value = "합의된 로맨스는 언제나 성공합니다"

$ws.Range("B21").Value = "Keyed"
$ws.Range("B22").Value = "Keyed"

$refFont = $ws.Range("B20").Font.Name
$ws.Range("B21").Font.Name = $refFont
$ws.Range("C21").Font.Name = $refFont
$ws.Range("E21").Font.Name = $refFont
$ws.Range("F21").Font.Name = $refFont
$ws.Range("B22").Font.Name = $refFont
$ws.Range("C22").Font.Name = $refFont
$ws.Range("E22").Font.Name = $refFont
$ws.Range("F22").Font.Name = $refFont

$ws.Range("F27").Select()
